$wb = $excel.ActiveWorkbook

# Duplicate the Slovakia sheet to serve as the template for the new Hungary sheet
$slovakia = $wb.Worksheets.Item("Slovakia")
$slovakia.Copy($null, $slovakia) | Out-Null

# The copy is inserted right after Slovakia and picks up a default name
$hungary = $wb.Worksheets.Item("Slovakia (2)")
$hungary.Name = "Hungary"

# Update the market name and NGC reference for Hungary
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Range("B4").Value = "NGC-4308/T3593/T3618/T3595/T3620"

# Remove the stray "LCD800" leftover row that Slovakia's sheet carried (row 10)
$hungary.Rows(10).Delete() | Out-Null

# B4 should carry the same bordered style as the other constant cells below it
$hungary.Range("B3").Copy() | Out-Null
$hungary.Range("B4").PasteSpecial(-4122) | Out-Null
$hungary.Application.CutCopyMode = $false

# Leave the new sheet's selection on B5, matching the authored workbook
$hungary.Range("B5").Select() | Out-Null

# Slovakia is no longer the active tab; its selection reverts to the whole sheet
$slovakia.Cells.Select() | Out-Null

# Hungary becomes the active sheet/tab
$hungary.Select() | Out-Null
$hungary.Range("B5").Select() | Out-Null
